# Add custom segments / Edit custom segments function for TTD
#
# The TTD sheet's sample/test row is replaced with a new "Add custom
# segment" style example row (Segment ID, Parent Segment ID == partner
# rate key, Segment Name/Description, and a Buyable boolean flag), and
# the old second example row is cleared out (only the formatted-but-empty
# Price cell remains). The previously-active sheet (Adobe AAM) loses focus
# and TTD becomes the selected/active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TTD")

# --- Row 3: new "Add Custom" example segment -----------------------------
$ws.Range("A3").Value = 20190401003
$ws.Range("B3").Value = "ttdratetest_partnerID_rate"
$ws.Range("C3").Value = "Test Segment 20190401003"
$ws.Range("D3").Value = "Test Segment 20190401003"
$ws.Range("E3").Value = $true

# Columns F/G/H are no longer used for this example row
$ws.Range("F3").ClearContents()
$ws.Range("G3").ClearContents()
$ws.Range("H3").ClearContents()

# --- Row 4: old second example row removed --------------------------------
# Keep G4's (empty) currency-formatted cell, clear everything else.
$ws.Range("A4:F4").ClearContents()
$ws.Range("G4").ClearContents()
$ws.Range("H4").ClearContents()

# --- Make TTD the active/selected sheet -----------------------------------
$ws.Activate() | Out-Null
$ws.Range("F9").Select() | Out-Null
